$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra helper columns F:I (G1/G2/G3/AGR headers and blank F2/F3 cells)
$ws.Range("F1:I3").EntireColumn.Delete()

# Replace the raw numeric revenue/net-income figures with their formatted
# dollar-string equivalents (stored as literal text, not numbers)
$ws.Range("B2").Value = "'$59,972,000"
$ws.Range("B2").Style = "Normal"

$ws.Range("B3").Value = "'-$373,705"
$ws.Range("B3").Style = "Normal"

$ws.Range("C2").Value = "'$76,033,000"
$ws.Range("C2").Style = "Normal"

$ws.Range("C3").Value = "'-$520,379"
$ws.Range("C3").Style = "Normal"

$ws.Range("D2").Value = "'$40,269,000"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'-$1,166,391"
$ws.Range("D3").Style = "Normal"

$ws.Range("E2").Value = "'$34,343,000"
$ws.Range("E2").Style = "Normal"

$ws.Range("E3").Value = "'-$579,646"
$ws.Range("E3").Style = "Normal"
